# Update column F (dSF) values on the active worksheet to match the
# re-pulled / re-computed data from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -1
    5  = 8
    7  = 1
    8  = 4
    14 = 6
    15 = -2
    16 = -1
    20 = 0
    23 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
